# Add a new "AUTH_06 - system error" test case row, and fix the
# Column C / Column D content swap (Sample Data vs Steps) for every
# existing test-case row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row right below AUTH_01 (row 2) for the new AUTH_06
#    "system error" scenario. This shifts the old rows 3-6 down to 4-7
#    and copies the formatting (incl. the green "PASS" style on col G)
#    from the row above, which is exactly what we want.
$ws.Rows.Item(3).Insert()

# 2) Rewrite every data row (2-7) with its final, correct values - note
#    column C now holds "Du Lieu Mau" (sample data) and column D holds
#    "Cac Buoc" (steps), matching the sheet's own header labels.

# Row 2 : AUTH_01 (unchanged scenario, columns C/D fixed)
$ws.Range("A2").Value = "AUTH_01"
$ws.Range("B2").Value = "Đăng nhập với input rỗng"
$ws.Range("C2").Value = 'User: "", Pass: "123"'
$ws.Range("D2").Value = "1. Nhập username rỗng`n2. Gọi hàm login"
$ws.Range("E2").Value = "Trạng thái: FAILED_CREDENTIALS"
$ws.Range("F2").Value = "OK"
$ws.Range("G2").Value = "PASS"

# Row 3 : AUTH_06 (new scenario - Database / system error)
$ws.Range("A3").Value = "AUTH_06"
$ws.Range("B3").Value = "Lỗi hệ thống (Database Error)"
$ws.Range("C3").Value = "User: any, Pass: any"
$ws.Range("D3").Value = "1. Mock DAO ném Exception`n2. Service catch và trả về lỗi hệ thống"
$ws.Range("E3").Value = "Trạng thái: FAILED_SYSTEM_ERROR"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"

# Row 4 : AUTH_05 (was row 3, columns C/D fixed)
$ws.Range("A4").Value = "AUTH_05"
$ws.Range("B4").Value = "Đăng nhập role lạ (Hacker/Lỗi data)"
$ws.Range("C4").Value = "User: hacker, Role: hacker_role"
$ws.Range("D4").Value = "1. Mock DAO trả về role 'hacker'`n2. Switch case check role"
$ws.Range("E4").Value = "Trạng thái: FAILED_INVALID_ROLE"
$ws.Range("F4").Value = "OK"
$ws.Range("G4").Value = "PASS"

# Row 5 : AUTH_02 (was row 4, columns C/D fixed)
$ws.Range("A5").Value = "AUTH_02"
$ws.Range("B5").Value = "Đăng nhập sai thông tin"
$ws.Range("C5").Value = "User: wrongUser, Pass: 123"
$ws.Range("D5").Value = "1. Mock DAO trả về null (không tìm thấy)`n2. Gọi login"
$ws.Range("E5").Value = "Trạng thái: FAILED_CREDENTIALS"
$ws.Range("F5").Value = "OK"
$ws.Range("G5").Value = "PASS"

# Row 6 : AUTH_04 (was row 5, columns C/D fixed)
$ws.Range("A6").Value = "AUTH_04"
$ws.Range("B6").Value = "Đăng nhập thành công (User thường)"
$ws.Range("C6").Value = "User: user1, Role: user"
$ws.Range("D6").Value = "1. Mock DAO trả về User thường`n2. Check Role"
$ws.Range("E6").Value = "Trạng thái: SUCCESS_USER"
$ws.Range("F6").Value = "OK"
$ws.Range("G6").Value = "PASS"

# Row 7 : AUTH_03 (was row 6, columns C/D fixed)
$ws.Range("A7").Value = "AUTH_03"
$ws.Range("B7").Value = "Đăng nhập thành công (Admin)"
$ws.Range("C7").Value = "User: admin, Role: admin"
$ws.Range("D7").Value = "1. Mock DAO trả về User Admin`n2. Check Role"
$ws.Range("E7").Value = "Trạng thái: SUCCESS_ADMIN"
$ws.Range("F7").Value = "OK"
$ws.Range("G7").Value = "PASS"

# 3) Re-fit the row heights back to the default after writing multi-line
#    text, so rows don't end up with a stray custom height.
$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(3).EntireRow.AutoFit()
$ws.Rows.Item(4).EntireRow.AutoFit()
$ws.Rows.Item(5).EntireRow.AutoFit()
$ws.Rows.Item(6).EntireRow.AutoFit()
$ws.Rows.Item(7).EntireRow.AutoFit()

# 4) Columns C and D effectively swapped meaning (C = sample data,
#    D = steps), so their best-fit widths swap too; column E's longest
#    value also changed slightly. Match the new best-fit widths as
#    closely as this engine's column-width rounding allows.
$ws.Columns.Item(3).ColumnWidth = 28.09
$ws.Columns.Item(4).ColumnWidth = 36.59
$ws.Columns.Item(5).ColumnWidth = 31.92

Write-Output "done"
